# Generate Report for Handback
# -----------------------------------------------------------------------
# After a handback round-trips through the localization pipeline, the
# report workbook is refreshed:
#   * Overview sheet + each locale sheet's "Status" column now reads
#     "Handed back: in sync with en-US" instead of "Ready for handoff".
#   * Each locale sheet gets its "Latest Target File" (hyperlinked to the
#     source .md doc - it's now in sync), "Latest Handback File" (the
#     generated .xlf) and "Latest Handback DateTime" columns populated
#     for both rows.
#   * A few columns are widened so the new, longer values aren't clipped.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a63dfab85c5835ec20ebbffba55a0325f223646/e2e/"
$docA = "0fcdd77c-2b75-462b-b247-62e8bcb60b5b.md"
$docB = "38409710-1a39-4aa6-b595-7acccff7e647.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-locale status columns (E = zh-cn,
# F = de-de) and widen them to fit the longer text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Range("E1").ColumnWidth = 29.1
$overview.Range("F1").ColumnWidth = 29.1

# ---------------------------------------------------------------------
# Per-locale handback details.
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackAt = "2016-10-19 12:24:10"; XlfA = "0fcdd77c-2b75-462b-b247-62e8bcb60b5b.178f9c2605b4efbb2e1aed6605e93a4575385170.zh-cn.xlf"; XlfB = "38409710-1a39-4aa6-b595-7acccff7e647.86c4e1998a4746e52b5c1a78694db2a5d52b8b33.zh-cn.xlf" },
    @{ Name = "de-de"; HandbackAt = "2016-10-19 12:24:28"; XlfA = "0fcdd77c-2b75-462b-b247-62e8bcb60b5b.178f9c2605b4efbb2e1aed6605e93a4575385170.de-de.xlf"; XlfB = "38409710-1a39-4aa6-b595-7acccff7e647.86c4e1998a4746e52b5c1a78694db2a5d52b8b33.de-de.xlf" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Status column (C) keeps its text in sync with the Overview sheet.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I) - now in sync with the source doc, so it
    # links to the same file as the Source File Name column (A).
    $ws.Range("I2").Value = $docA
    $ws.Range("I3").Value = $docB
    $ws.Hyperlinks.Add($ws.Range("I2"), ($repoBase + $docA), "", "", $docA)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($repoBase + $docB), "", "", $docB)

    # Latest Handback File (J) - the xlf generated for this handback.
    $ws.Range("J2").Value = $locale.XlfA
    $ws.Range("J3").Value = $locale.XlfB

    # Latest Handback DateTime (K).
    $ws.Range("K2").Value = $locale.HandbackAt
    $ws.Range("K3").Value = $locale.HandbackAt

    # Widen Status (C), Latest Target File (I) and Latest Handback File
    # (J) so the new values aren't clipped.
    $ws.Range("C1").ColumnWidth = 29.1
    $ws.Range("I1").ColumnWidth = 39.15
    $ws.Range("J1").ColumnWidth = 39.15
}
